$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 410; this shifts existing rows 410..438 down to 411..439
$ws.Rows.Item(410).Insert()

# Populate the newly inserted row 410 with the new weekly data point
$ws.Cells.Item(410, 1).Value = 4
$ws.Cells.Item(410, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(410, 3).Value = "Los Lagos"
$ws.Cells.Item(410, 4).Value = "2023-04-05"
$ws.Cells.Item(410, 5).Value = 10
$ws.Cells.Item(410, 6).Value = 100112017
$ws.Cells.Item(410, 7).Value = "Apio"
$ws.Cells.Item(410, 8).Value = "Americana (o)"
$ws.Cells.Item(410, 9).Value = "Primera"
$ws.Cells.Item(410, 10).Value = 15
$ws.Cells.Item(410, 11).Value = 12000
$ws.Cells.Item(410, 12).Value = 12000
$ws.Cells.Item(410, 13).Value = 12000
$ws.Cells.Item(410, 14).Value = "`$/docena de matas"
$ws.Cells.Item(410, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(410, 16).Value = 2000
$ws.Cells.Item(410, 17).Value = 6
$ws.Cells.Item(410, 18).Value = "Hortaliza"
